# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns, matching the style of the existing
# header row (bold, bordered, centered) by copying the formatting from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Numeric data for rows 2-28.
$data = @{
    2  = @(1, 5)
    3  = @(1, 6)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(4, 7)
    7  = @(1, 4)
    8  = @(5, 7)
    9  = @(5, 6)
    10 = @(8, 8)
    11 = @(7, 9)
    12 = @(6, 6)
    13 = @(6, 6)
    14 = @(6, 7)
    15 = @(6, 7)
    16 = @(5, 6)
    17 = @(6, 8)
    18 = @(11, 12)
    19 = @(4, 7)
    20 = @(5, 7)
    21 = @(8, 9)
    22 = @(8, 9)
    23 = @(5, 6)
    24 = @(5, 7)
    25 = @(1, 4)
    26 = @(1, 5)
    27 = @(1, 4)
    28 = @(1, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}

Write-Output "done"
